$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "69.960.16" }
    @{ Cell = "E2"; Value = "  -0.65%  " }
    @{ Cell = "D3"; Value = "3.544.30" }
    @{ Cell = "E3"; Value = "  -0.82%  " }
    @{ Cell = "E4"; Value = "  +0.03%  " }
    @{ Cell = "D5"; Value = "605.93" }
    @{ Cell = "E5"; Value = "  +2.60%  " }
    @{ Cell = "D6"; Value = "184.73" }
    @{ Cell = "E6"; Value = "  -1.38%  " }
    @{ Cell = "D7"; Value = "3.543.16" }
    @{ Cell = "E7"; Value = "  -0.44%  " }
    @{ Cell = "D8"; Value = "0.615" }
    @{ Cell = "E8"; Value = "  -1.02%  " }
    @{ Cell = "E9"; Value = "  -0.01%  " }
    @{ Cell = "D10"; Value = "0.214" }
    @{ Cell = "E10"; Value = "  +7.18%  " }
    @{ Cell = "D11"; Value = "0.640" }
    @{ Cell = "E11"; Value = "  -1.19%  " }
    @{ Cell = "D12"; Value = "53.52" }
    @{ Cell = "E12"; Value = "  -2.44%  " }
    @{ Cell = "D13"; Value = "0.0000308" }
    @{ Cell = "E13"; Value = "  +0.02%  " }
    @{ Cell = "D14"; Value = "9.43" }
    @{ Cell = "E14"; Value = "  -1.63%  " }
    @{ Cell = "D15"; Value = "4.120.03" }
    @{ Cell = "E15"; Value = "  -0.41%  " }
    @{ Cell = "D16"; Value = "70.089.93" }
    @{ Cell = "E16"; Value = "  -0.38%  " }
    @{ Cell = "D17"; Value = "3.564.86" }
    @{ Cell = "E17"; Value = "  -0.29%  " }
    @{ Cell = "B18"; Value = "Chainlink" }
    @{ Cell = "C18"; Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link" }
    @{ Cell = "D18"; Value = "18.88" }
    @{ Cell = "E18"; Value = "  -3.10%  " }
    @{ Cell = "B19"; Value = "Uniswap" }
    @{ Cell = "C19"; Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni" }
    @{ Cell = "D19"; Value = "12.60" }
    @{ Cell = "E19"; Value = "  +0.68%  " }
    @{ Cell = "D20"; Value = "577.11" }
    @{ Cell = "E20"; Value = "  +5.54%  " }
    @{ Cell = "E21"; Value = "  +0.17%  " }
    @{ Cell = "D22"; Value = "0.987" }
    @{ Cell = "E22"; Value = "  -3.42%  " }
    @{ Cell = "D23"; Value = "17.33" }
    @{ Cell = "E23"; Value = "  -3.86%  " }
    @{ Cell = "D24"; Value = "4.68" }
    @{ Cell = "E24"; Value = "  +0.16%  " }
    @{ Cell = "D25"; Value = "4.82" }
    @{ Cell = "E25"; Value = "  -2.25%  " }
    @{ Cell = "D26"; Value = "94.39" }
    @{ Cell = "E26"; Value = "  -1.79%  " }
    @{ Cell = "D27"; Value = "2.93" }
    @{ Cell = "E27"; Value = "  -2.54%  " }
    @{ Cell = "D28"; Value = "10.94" }
    @{ Cell = "E28"; Value = "  -4.82%  " }
    @{ Cell = "D29"; Value = "9.32" }
    @{ Cell = "E29"; Value = "  +1.66%  " }
    @{ Cell = "D30"; Value = "32.04" }
    @{ Cell = "E30"; Value = "  -0.69%  " }
    @{ Cell = "D31"; Value = "7.00" }
    @{ Cell = "E31"; Value = "  -4.76%  " }
    @{ Cell = "D32"; Value = "12.18" }
    @{ Cell = "E32"; Value = "  -3.09%  " }
    @{ Cell = "D33"; Value = "0.113" }
    @{ Cell = "E33"; Value = "  -1.19%  " }
    @{ Cell = "D34"; Value = "63.27" }
    @{ Cell = "E34"; Value = "  -3.01%  " }
    @{ Cell = "D35"; Value = "3.66" }
    @{ Cell = "E35"; Value = "  +18.27%  " }
    @{ Cell = "E36"; Value = "  +0.12%  " }
    @{ Cell = "D37"; Value = "531.09" }
    @{ Cell = "E37"; Value = "  -3.60%  " }
    @{ Cell = "D38"; Value = "0.400" }
    @{ Cell = "E38"; Value = "  -4.35%  " }
    @{ Cell = "D39"; Value = "0.999" }
    @{ Cell = "E39"; Value = "  -0.10%  " }
    @{ Cell = "D40"; Value = "36.85" }
    @{ Cell = "E40"; Value = "  -4.14%  " }
    @{ Cell = "D41"; Value = "0.0₃0780" }
    @{ Cell = "E41"; Value = "  +1.12%  " }
    @{ Cell = "D42"; Value = "3.529.93" }
    @{ Cell = "E42"; Value = "  +4.72%  " }
    @{ Cell = "D43"; Value = "3.50" }
    @{ Cell = "E43"; Value = "  +3.33%  " }
    @{ Cell = "D44"; Value = "0.135" }
    @{ Cell = "E44"; Value = "  -0.06%  " }
    @{ Cell = "D45"; Value = "0.0456" }
    @{ Cell = "E45"; Value = "  +1.45%  " }
    @{ Cell = "D46"; Value = "3.45" }
    @{ Cell = "E46"; Value = "  -3.79%  " }
    @{ Cell = "D47"; Value = "2.90" }
    @{ Cell = "E47"; Value = "  -3.35%  " }
    @{ Cell = "E48"; Value = "  +2.22%  " }
    @{ Cell = "D49"; Value = "9.14" }
    @{ Cell = "E49"; Value = "  -0.89%  " }
    @{ Cell = "E50"; Value = "  +0.46%  " }
    @{ Cell = "B51"; Value = "OceanProtocol" }
    @{ Cell = "C51"; Value = "https://coinranking.com/coin/aAKLSV5-0+oceanprotocol-ocean" }
    @{ Cell = "D51"; Value = "1.40" }
    @{ Cell = "E51"; Value = "  -4.60%  " }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
    $cell.Style = $origStyle
}
